# SOlo se puede ver sus proyectos
# Adds a new "projects.listall" permission row to the Permissions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New permission row (row 67 -> id 66) ---------------------------------
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "projects.listall"
$ws.Cells.Item(67, 3).Value = "El usuario podrá listar todos los proyectos. Si esta desactivado, solo puede ver sus asignados."

# Column C got noticeably wider once this longer description was added
# (Excel's "best fit" autosize kicking in for the new longest string).
$ws.Columns.Item(3).ColumnWidth = 83.6

# --- View state: scroll position & selection moved while editing ---------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C64").Select()
